# DEV-3439: remove the standalone "IntervalValue"/hasInterval row and the
# "Geonames" gui_element row, collapsing the data block by one row.
#
# Original sheet (rows 1-5):
#   1 header
#   2 hasBoolean / Checkbox
#   3 hasGeoname / Geoname link ... gui_element=Geonames
#   4 hasInterval / Time interval ... gui_element=Geonames   <- removed entirely
#   5 (blank row holding only the trailing "  " placeholder strings)  <- becomes new row 4
#
# Target sheet (rows 1-4):
#   1 header
#   2 hasBoolean / Checkbox                              (unchanged values)
#   3 hasGeoname / Geoname link ... gui_element=Checkbox  (gui_element fixed)
#   4 (blank placeholder row, shifted up from old row 5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "hasInterval" row. This removes row 4 and shifts the old
# row 5 (the blank "  " placeholder row) up to become the new row 4, matching
# the target dimension A1:O4.
$ws.Rows.Item(4).Delete()

# The gui_element value in row 3 (N3) incorrectly pointed at "Geonames";
# it should reference the "Checkbox" gui element, same as row 2.
$ws.Range("N3").Value = "Checkbox"

# J2/J3 were using a stray duplicate font/style (direct formatting only used
# by this column in these two cells); normalise them to the same format used
# by the rest of column J / the sibling "hasValue" cells (L2/L3), which also
# lets the now-unused font drop out of the style table.
$ws.Range("J2").Font.Color = 0
$ws.Range("J3").Font.Color = 0
